# Append a fresh scrape pass (2026-01-16 12:40:04) to the "ランサーズ" sheet.
# The scraper re-writes the whole listing in descending-priority order, so
# several brand-new postings are interleaved with the previously-seen ones
# (which simply shift down). We rebuild rows 2-13 in full and re-create the
# URL hyperlinks to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamp = "2026-01-16 12:40:04"

# Columns: Title, Category, Price, Deadline, URL, Score, Skills(optional)
$rows = @(
    @("製造業向け図面自動生成システムの開発・ツール化を支援してくださるエンジニア募集(AI/バックエンド)", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5460562", 435, "🔥AI,Ai ◆ツール,開発"),
    @("【募集】Python / Docker 日次データ スクレイピングシステム構築", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5469627", 248, "🔥Python ◆スクレイピング"),
    @("【募集】RPAツール「RoboTANGO」設定代行の専門家を探しています", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5405023", 178, "★bot ◆ツール"),
    @("施設管理・現場業務向け チェックリスト業務の自動化・報告書作成システム開発エンジニア募集", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5460563", 220, "◆開発,システム開発 ◇管理"),
    @("【Windows/Wacom】署名画像から筆順解析図を作成する業務用アプリ開発", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5472804", 100, "◆開発 ◇アプリ"),
    @("署名画像から筆順を可視化するアプリ開発者募集", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5472080", 100, "◆開発 ◇アプリ"),
    @("スマホカラオケ予約Webアプリ開発のフリーランス募集(使用するのは個人の集まりで趣味で使う程度です)", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5472431", 88, "◆開発 ◇アプリ"),
    @("【自動化】申込書AからBへの転写をエクセルで実現したい", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5473042", 83, "◆自動化"),
    @("初回 Webサーバ管理エンジニア", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5472544", 45, "◇管理"),
    @("【急募】Flutterflowの扱えるノーコードエンジニアを探しています!", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5472976", 25, ""),
    @("m.2 SSD基板の設計", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5472120", 25, ""),
    @("《長期レギュラー》公的機関Web運用の要となる、ディレクター募集", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5472958", 18, "")
)

# Drop the old hyperlink objects; we'll recreate one per row below so the
# relationship ids line up with the rebuilt F column.
$ws.Hyperlinks.Delete()

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $timestamp
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]

    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $row[4])
    $ws.Cells.Item($r, 6).Style = "Hyperlink"

    $ws.Cells.Item($r, 7).Value = $row[5]

    if ($row[6] -ne "") {
        $ws.Cells.Item($r, 8).Value = $row[6]
    }

    $r = $r + 1
}
